$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

$shp = $s.Shapes.AddTextbox(1, 711.6923622047244, 14.209370078740157, 106.15377952755905, 29.081259842519685)
$shp.TextFrame.TextRange.Text = "Scale-up?"
$shp.TextFrame.TextRange.LanguageID = "en-GB"
$shp.Fill.Visible = $false
$shp.TextFrame.AutoSize = 1
$shp.TextFrame.WordWrap = -1
